# Rename existing sheet strategy_id-5008 -> strategy_id-5007,
# then add a new sheet strategy_id-5009 (copy of strategy_id-5007) right after it.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("strategy_id-5008")
$src.Name = "strategy_id-5007"

# Duplicate the sheet (copy placed after itself) to create the new strategy_id-5009 tab.
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item($src.Index + 1)
$newSheet.Name = "strategy_id-5009"
